$wb = $excel.ActiveWorkbook
$wsTables = $wb.Worksheets.Item("Tables")
$wsProcedures = $wb.Worksheets.Item("Procedures")

# --- Procedures sheet (sheet2): append 12 new stored-procedure names ----------
$procNames = @(
    "FOX_PROC_GET_CASE_TRATEMENT_TEAM",
    "FOX_PROC_GET_CONSENT_TO_CARE_DETAILS_BY_CASE_ID",
    "FOX_PROC_GET_CONSENT_TO_CARE_DOCUMENTS_INFO",
    "FOX_PROC_GET_CONSENT_TO_CARE_INFO_BY_CASE_ID",
    "FOX_PROC_GET_CONSENT_TO_CARE_INFO_BY_CASE_ID_AND_SEND_TO",
    "FOX_PROC_GET_CONSENT_TO_CARE_INFO_BY_CONSENT_TO_CARE_ID",
    "FOX_PROC_GET_INSURANCE_DETAILS_FOR_CONSENT_TO_CARE",
    "FOX_PROC_GET_PATINET_CONTACT_DETAILS",
    "FOX_PROC_GET_SERVICE_CONFIGURATION_CONSENT_TO_CARE",
    "FOX_PROC_GET_USER_ID_BY_PROVIDER_CODE",
    "FOX_PROC_INSERT_CONSENT_TO_CARE_TASK",
    "FOX_PROC_UPDATE_TASK_LOG"
)

$row = 478
foreach ($name in $procNames) {
    $wsProcedures.Cells.Item($row, 1).Value = $name
    $row = $row + 1
}

# --- Tables sheet (sheet1): append 3 new table names --------------------------
$tableNames = @(
    "FOX_TBL_CONSENT_TO_CARE",
    "FOX_TBL_CONSENT_TO_CARE_DOCUMENTS",
    "FOX_TBL_CONSENT_TO_CARE_STATUS"
)

$row = 215
foreach ($name in $tableNames) {
    $wsTables.Cells.Item($row, 1).Value = $name
    $row = $row + 1
}

# --- View / selection state ----------------------------------------------------
# Procedures tab is no longer the active tab; its selection moves to the new last row.
$wsProcedures.Activate()
$wsProcedures.Range("A489").Select()

# Tables tab becomes the active tab, selection on the new last row.
$wsTables.Activate()
$wsTables.Range("A217").Select()
